# daily auto push: 2026-02-03 03:12 UTC
# Insert two new daily records (2026/02/02 23:00 and 2026/02/03 07:00)
# right before the existing 2026/12/29 block, pushing all later rows down
# by two and extending the sheet from 786 to 788 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertAt = 745

# Use "copy an existing row, then Insert()" so the new rows inherit the
# same cell typing (plain text date/weekday, not an auto-converted date
# serial) as the rest of the sheet, and everything from $insertAt onward
# shifts down by one row each time.
$ws.Rows.Item($insertAt).Copy()
$ws.Rows.Item($insertAt).Insert()

$ws.Rows.Item($insertAt + 1).Copy()
$ws.Rows.Item($insertAt + 1).Insert()

# Force columns A/B of the two new rows to Text so assigning date-like or
# weekday-kanji strings doesn't get reinterpreted as a date value.
$newRows = $ws.Range("A" + $insertAt + ":B" + ($insertAt + 1))
$newRows.NumberFormat = "@"

$ws.Cells.Item($insertAt, 1).Value = "2026/02/02"
$ws.Cells.Item($insertAt, 2).Value = "月"
$ws.Cells.Item($insertAt, 3).Value = 23
$ws.Cells.Item($insertAt, 4).Value = 20

$ws.Cells.Item($insertAt + 1, 1).Value = "2026/02/03"
$ws.Cells.Item($insertAt + 1, 2).Value = "火"
$ws.Cells.Item($insertAt + 1, 3).Value = 7
$ws.Cells.Item($insertAt + 1, 4).Value = 22

# Drop the temporary Text number format again so the new cells end up with
# the same (default) style as the rest of the data rows.
$newRows.ClearFormats()
